$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price / volume updates ---
$ws.Range("D2").Value = "29.341.46"
$ws.Range("E2").Value = "  -0.81%  "
$ws.Range("D3").Value = "1.839.55"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.92"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6249"
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07375"
$ws.Range("E8").Value = "  -1.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2883"
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.75"
$ws.Range("E10").Value = "  -0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07720"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "1.831.43"
$ws.Range("E12").Value = "  -0.50%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.953"
$ws.Range("E13").Value = "  -1.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6663"
$ws.Range("E14").Value = "  -2.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001035"
$ws.Range("E15").Value = "  -0.89%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.44"
$ws.Range("E16").Value = "  -1.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.231"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "29.310.80"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "234.10"
$ws.Range("E19").Value = "  +1.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.26"
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.280"
$ws.Range("E22").Value = "  -3.95%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "157.22"
$ws.Range("E24").Value = "  -1.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.449"
$ws.Range("E25").Value = "  -1.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1334"
$ws.Range("E26").Value = "  -3.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.27"
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07149"
$ws.Range("E28").Value = "  +8.11%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.487"
$ws.Range("E29").Value = "  +2.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.479"
$ws.Range("E30").Value = "  -0.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.153"
$ws.Range("E33").Value = "  +0.19%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.812"
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7043"
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.587"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01826"
$ws.Range("E37").Value = "  -2.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.783"
$ws.Range("E38").Value = "  -2.03%  "
$ws.Range("D39").Value = "1.231.16"
$ws.Range("E39").Value = "  -3.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.749"
$ws.Range("E40").Value = "  -1.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9465"
$ws.Range("E41").Value = "  +1.09%  "
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.02"
$ws.Range("E45").Value = "  -2.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000118"
$ws.Range("E46").Value = "  +2.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.923"
$ws.Range("E47").Value = "  -2.82%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.684"
$ws.Range("E48").Value = "  -3.30%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.864"
$ws.Range("E49").Value = "  -2.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1126"
$ws.Range("E50").Value = "  -3.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3870"
$ws.Range("E51").Value = "  -2.48%  "

# --- Row 31 / 32 swap: Filecoin <-> InternetComputer(DFINITY) ---
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.030"
$ws.Range("E31").Value = "  -1.96%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.022"
$ws.Range("E32").Value = "  -2.18%  "

# --- Row 43 / 44 swap: RocketPoolETH <-> Quant ---
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.90"
$ws.Range("E43").Value = "  -0.36%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.978.52"
$ws.Range("E44").Value = "  -1.93%  "
